$wb = $excel.ActiveWorkbook

# --- Update CaseDetailStat sheet (F2): file size now shown with unit suffix ---
$wsStat = $wb.Worksheets.Item("CaseDetailStat")
$wsStat.Range("F2").Value = "105.75 KB"

# --- Update CaseDetailStat_Message sheet (A28): Cypher query text updated to
#     compute a human readable size (value + unit) instead of a raw byte count ---
$wsMsg = $wb.Worksheets.Item("CaseDetailStat_Message")

$newQuery = @'
MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent)
WHERE c.case_id IN ['NCATS-COP01CCB050022']
WITH
['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
toInteger(floor(log(f.file_size)/log(1024))) as i,
2 as precision,
f,parent
WITH f.file_size /(1024^i) AS value, 10^precision AS factor, units[i] as unit,f,parent
RETURN f.file_name AS `File Name` ,f.file_type AS `File Type`,head(labels(parent)) AS `Association`, f.file_description AS `Description`,f.file_format AS Format,round(factor * value)/factor+ +unit AS Size
'@

$wsMsg.Range("A28").Value = $newQuery
